$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "sitemap,xml et robots,txt" -> "robots,txt", and add G9 = "x"
$ws.Range("B9").Value = "robots,txt"
$ws.Range("G9").Value = "x"

# Remove the two obsolete rows ("keyword à 1px / black hat SEO" and
# "déplacer les scripts JS"), which shifts rows 15-20 up to 13-18.
$ws.Rows("13:14").Delete()

# Insert two fresh rows for the new audit items, pushing old row 18
# ("(SEO ou accessiblité ?)") back down to row 18.
$ws.Rows("16:17").Insert()
$ws.Rows("16:17").RowHeight = 15.75

$ws.Range("A16").Value = "(SEO ou accessiblité ?)"
$ws.Range("B16").Value = "minifier le css"

$ws.Range("A17").Value = "(SEO ou accessiblité ?)"
$ws.Range("B17").Value = "Bug dans bloc.js"

# Old trailing placeholder rows 19/20 are now blank.
$ws.Range("A19").ClearContents()
$ws.Range("A20").ClearContents()

# Selection moved to B18 in the saved file.
$ws.Range("B18").Select()
